$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

$win.FreezePanes = $false
$ws.Range("B2").Select()
$win.FreezePanes = $true
Write-Host "post freeze: Split=$($win.SplitRow),$($win.SplitColumn) Scroll=$($win.ScrollRow),$($win.ScrollColumn)"

# scroll down+right by 5 rows, 3 cols using LargeScroll/SmallScroll
$win.SmallScroll(5,0,3,0) | Out-Null
Write-Host "post smallscroll: Scroll=$($win.ScrollRow),$($win.ScrollColumn)"
